$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Change shared string "E7760" -> "E7420" wherever it appears (column H, rows 2-27)
for ($r = 2; $r -le 27; $r++) {
    $cell = $ws.Cells.Item($r, 8)  # column H
    if ($cell.Value2 -eq "E7760") {
        $cell.Value = "E7420"
    }
}

# 2. Update the selection to H2:H27 with active cell H2
$ws.Range("H2:H27").Select()
